$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Overwrite row 4 with the data that used to be in row 5
$ws.Range("A4").Value = "even_MAG-GUT81714.fa"
$ws.Range("B4").Value = 0.638197623694887
$ws.Range("C4").Value = 0.03776221813686578
$ws.Range("D4").Value = 0.001305826889599017
$ws.Range("E4").Value = 0.008040952217973926
$ws.Range("F4").Value = 0.003110745400189174
$ws.Range("G4").Value = 0.1000081523671545
$ws.Range("H4").Value = 0.2115744812933307
$ws.Range("I4").Value = 0.638197623694887
$ws.Range("J4").Value = "s__Agathobaculum butyriciproducens"
$ws.Range("K4").Value = "s__Agathobaculum butyriciproducens"

# Remove the now-duplicate row 5 entirely
$ws.Rows.Item(5).Delete()
